$d = $word.ActiveDocument

# Collapsed range positioned right at the very end of the document body
# (before the final paragraph mark), so the inserted XML is appended
# after the last existing paragraph and before the sectPr.
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:cs="B Titr"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">تاریخ: </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>28</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> آبان 1401</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:cs="B Titr"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr" w:hint="cs"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:cs="B Titr"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$endRange.InsertXML($xml)
